$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row-level updates taken from the "cryptos list" refresh.
# D = new price text (or $null if unchanged), E = new 1h volume/change text.
$updates = @(
    @{ Row = 2; D = "46.961.47"; E = "  +5.07%  " },
    @{ Row = 3; D = "2.355.77"; E = "  +4.66%  " },
    @{ Row = 4; D = $null; E = "  -0.95%  " },
    @{ Row = 5; D = "307.89"; E = "  +0.41%  " },
    @{ Row = 6; D = "99.18"; E = "  +4.38%  " },
    @{ Row = 7; D = $null; E = "  +1.55%  " },
    @{ Row = 8; D = $null; E = "  -0.60%  " },
    @{ Row = 9; D = "0.540"; E = "  +4.34%  " },
    @{ Row = 10; D = "36.19"; E = "  +3.51%  " },
    @{ Row = 11; D = "0.0811"; E = "  +1.02%  " },
    @{ Row = 12; D = $null; E = "  +3.53%  " },
    @{ Row = 13; D = $null; E = "  -0.14%  " },
    @{ Row = 14; D = "2.713.79"; E = "  +4.58%  " },
    @{ Row = 15; D = "2.353.74"; E = "  -1.69%  " },
    @{ Row = 16; D = "14.39"; E = "  +5.79%  " },
    @{ Row = 17; D = $null; E = "  +0.10%  " },
    @{ Row = 18; D = "46.842.31"; E = "  +5.32%  " },
    @{ Row = 19; D = "13.59"; E = "  +15.80%  " },
    @{ Row = 20; D = $null; E = "  +1.72%  " },
    @{ Row = 21; D = $null; E = "  +0.27%  " },
    @{ Row = 22; D = "67.08"; E = "  +2.58%  " },
    @{ Row = 23; D = "245.80"; E = "  +3.31%  " },
    @{ Row = 24; D = "3.00"; E = "  +1.50%  " },
    @{ Row = 25; D = "2.01"; E = "  +1.75%  " },
    @{ Row = 26; D = $null; E = "  -0.25%  " },
    @{ Row = 27; D = "42.27"; E = "  +13.70%  " },
    @{ Row = 29; D = $null; E = "  +1.64%  " },
    @{ Row = 30; D = "20.27"; E = "  +1.37%  " },
    @{ Row = 31; D = "5.79"; E = "  -2.15%  " },
    @{ Row = 32; D = "152.24"; E = "  +2.49%  " },
    @{ Row = 33; D = "0.0819"; E = "  +4.28%  " },
    @{ Row = 34; D = "2.62"; E = "  +0.03%  " },
    @{ Row = 35; D = "3.17"; E = "  -1.84%  " },
    @{ Row = 36; D = "0.112"; E = "  +2.89%  " },
    @{ Row = 37; D = $null; E = "  +0.61%  " },
    @{ Row = 38; D = $null; E = "  -0.81%  " },
    @{ Row = 39; D = "4.12"; E = "  +8.75%  " },
    @{ Row = 40; D = $null; E = "  +6.46%  " },
    @{ Row = 41; D = $null; E = "  +2.99%  " },
    @{ Row = 42; D = "14.11"; E = "  -7.84%  " },
    @{ Row = 43; D = "0.998"; E = "  -0.83%  " },
    @{ Row = 44; D = "1.873.95"; E = "  +3.86%  " },
    @{ Row = 45; D = $null; E = "  +11.49%  " },
    @{ Row = 46; D = "0.200"; E = "  +6.37%  " },
    @{ Row = 47; D = "81.24"; E = "  -0.86%  " },
    @{ Row = 48; D = "74.05"; E = "  +7.76%  " },
    @{ Row = 49; D = $null; E = "  +2.43%  " },
    @{ Row = 50; D = "99.58"; E = "  +0.92%  " },
    @{ Row = 51; D = "55.82"; E = "  +3.34%  " }
)

foreach ($u in $updates) {
    $r = $u.Row

    if ($null -ne $u.D) {
        $text = $u.D
        # These "Price" cells are stored as plain text (e.g. "307.89", "0.0811").
        # If the replacement text looks like a plain number, Excel's
        # auto-detection would otherwise silently convert the cell to a real
        # number, so force text entry the same way a leading apostrophe would
        # in the UI.
        if ($text -match '^[0-9]+(\.[0-9]+)?$') {
            $ws.Range("D$r").Value = "'" + $text
        } else {
            $ws.Range("D$r").Value = $text
        }
    }

    $ws.Range("E$r").Value = $u.E
}
